$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 435.7143
$ws.Range("I4").Value = 453
$ws.Range("K4").Value = 453
$ws.Range("M4").Value = -339
$ws.Range("H33").Value = 1042.8667
$ws.Range("I33").Value = 520.1429000000001
$ws.Range("J33").Value = 1500.25
$ws.Range("K33").Value = 520.1429000000001
$ws.Range("L33").Value = 1500.25
$ws.Range("M33").Value = -291.1429000000001
$ws.Range("N33").Value = -1958.25
$ws.Range("H39").Value = 17.666666
$ws.Range("I39").Value = 17.666666
$ws.Range("K39").Value = 52.999998
$ws.Range("M39").Value = 243.000002
$ws.Range("H40").Value = 3839.1
$ws.Range("I40").Value = 3133.3333
$ws.Range("K40").Value = 3133.3333
$ws.Range("M40").Value = -2958.3333
$ws.Range("H98").Value = 2904.4707
$ws.Range("I98").Value = 3024.75
$ws.Range("K98").Value = 3024.75
$ws.Range("M98").Value = -1526.75
$ws.Range("H101").Value = 409.5
$ws.Range("I101").Value = 409.5
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1228.5
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H122").Value = 2904.4707
$ws.Range("I122").Value = 3024.75
$ws.Range("K122").Value = 9074.25
$ws.Range("M122").Value = -6624.25
$ws.Range("H132").Value = 5743793
$ws.Range("I132").Value = 7557246
$ws.Range("K132").Value = 22671738
$ws.Range("M132").Value = -22669208

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 2328.4
$ws.Range("I31").Value = 2328.4
$ws.Range("K31").Value = 2328.4
$ws.Range("M31").Value = -2034.4
$ws.Range("H32").Value = 24521.355
$ws.Range("I32").Value = 26688.902
$ws.Range("K32").Value = 26688.902
$ws.Range("M32").Value = -26401.902
$ws.Range("H74").Value = 256139.3
$ws.Range("I74").Value = 401100.34
$ws.Range("J74").Value = 14537.556
$ws.Range("K74").Value = 401100.34
$ws.Range("L74").Value = 14537.556
$ws.Range("M74").Value = -400226.34
$ws.Range("N74").Value = -16285.556
$ws.Range("H77").Value = 256139.3
$ws.Range("I77").Value = 401100.34
$ws.Range("J77").Value = 14537.556
$ws.Range("K77").Value = 2005501.7
$ws.Range("L77").Value = 72687.78
$ws.Range("M77").Value = -2001133.7
$ws.Range("N77").Value = -81423.78
$ws.Range("H122").Value = 2789
$ws.Range("J122").Value = 2638.5
$ws.Range("L122").Value = 7915.5
$ws.Range("N122").Value = -12815.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 178099
$ws.Range("J59").Value = 178099
$ws.Range("L59").Value = 178099
$ws.Range("N59").Value = -179793
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H94").Value = 1938.5
$ws.Range("I94").Value = 1938.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1938.5
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H102").Value = 27499.5
$ws.Range("I102").Value = 27499.5
$ws.Range("K102").Value = 27499.5
$ws.Range("M102").Value = -24254.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 5666
$ws.Range("J26").Value = 5666
$ws.Range("L26").Value = 5666
$ws.Range("N26").Value = -6240
$ws.Range("H58").Value = 15965.392
$ws.Range("J58").Value = 82978
$ws.Range("L58").Value = 82978
$ws.Range("N58").Value = -83384
$ws.Range("H132").Value = 145031.58
$ws.Range("I132").Value = 201044.4
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 603133.2
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -600603.2
$ws.Range("N132").Value = -20058.5
$ws.Range("H136").Value = 15965.392
$ws.Range("J136").Value = 82978
$ws.Range("L136").Value = 248934
$ws.Range("N136").Value = -254034

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1983.3334
$ws.Range("I55").Value = 975
$ws.Range("K55").Value = 2925
$ws.Range("M55").Value = -2748
$ws.Range("H92").Value = 1612.7142
$ws.Range("J92").Value = 1158.75
$ws.Range("L92").Value = 3476.25
$ws.Range("N92").Value = -5972.25
$ws.Range("H115").Value = 1077
$ws.Range("I115").Value = 200
$ws.Range("J115").Value = 1515.5
$ws.Range("K115").Value = 600
$ws.Range("L115").Value = 4546.5
$ws.Range("M115").Value = 575
$ws.Range("N115").Value = -6896.5
$ws.Range("H123").Value = 1805.3334
$ws.Range("I123").Value = 583
$ws.Range("K123").Value = 1749
$ws.Range("M123").Value = 701
$ws.Range("H132").Value = 2744.6667
$ws.Range("J132").Value = 1292.8
$ws.Range("L132").Value = 11635.2
$ws.Range("N132").Value = -16695.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 33166.5
$ws.Range("I80").Value = 1500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -502
$ws.Range("H83").Value = 33166.5
$ws.Range("I83").Value = 1500
$ws.Range("K83").Value = 7500
$ws.Range("M83").Value = -2508
$ws.Range("H102").Value = 31065.477
$ws.Range("J102").Value = 3172.5
$ws.Range("L102").Value = 3172.5
$ws.Range("N102").Value = -6416.5
$ws.Range("H132").Value = 4903.077
$ws.Range("I132").Value = 4885.4546
$ws.Range("K132").Value = 14656.3638
$ws.Range("M132").Value = -12126.3638

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2735.85
$ws.Range("I46").Value = 800.2727
$ws.Range("J46").Value = 5101.5557
$ws.Range("K46").Value = 800.2727
$ws.Range("L46").Value = 5101.5557
$ws.Range("M46").Value = -612.2727
$ws.Range("N46").Value = -5477.5557
$ws.Range("H55").Value = 871.5789
$ws.Range("I55").Value = 241.11111
$ws.Range("K55").Value = 241.11111
$ws.Range("M55").Value = -68.11111
$ws.Range("H119").Value = 101500
$ws.Range("J119").Value = 101500
$ws.Range("L119").Value = 101500
$ws.Range("N119").Value = -111176

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 97266.92999999999
$ws.Range("I122").Value = 166482.25
$ws.Range("J122").Value = 4979.8335
$ws.Range("K122").Value = 499446.75
$ws.Range("L122").Value = 14939.5005
$ws.Range("M122").Value = -496996.75
$ws.Range("N122").Value = -19839.5005
$ws.Range("H126").Value = 181349
$ws.Range("I126").Value = 1736.7778
$ws.Range("J126").Value = 504651
$ws.Range("K126").Value = 5210.3334
$ws.Range("L126").Value = 1513953
$ws.Range("M126").Value = -2740.3334
$ws.Range("N126").Value = -1518893
$ws.Range("H132").Value = 1859.4642
$ws.Range("I132").Value = 928.55
$ws.Range("K132").Value = 2785.65
$ws.Range("M132").Value = -255.6499999999996
$ws.Range("H136").Value = 14293.955
$ws.Range("I136").Value = 14721.6045
$ws.Range("J136").Value = 5099.5
$ws.Range("K136").Value = 44164.8135
$ws.Range("L136").Value = 15298.5
$ws.Range("M136").Value = -41614.8135
$ws.Range("N136").Value = -20398.5
